$wb = $excel.ActiveWorkbook

$ws_sheet1 = $wb.Worksheets.Item("展览")
$ws_sheet1.Range("F3").Value = 257
$ws_sheet1.Range("F4").Value = 81
$ws_sheet1.Range("F5").Value = 9626
$ws_sheet1.Range("F6").Value = 642
$ws_sheet1.Range("C7").Value = "北京·动画电影《钢管公主》专场活动"
$ws_sheet1.Range("D7").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws_sheet1.Range("E7").Value = "2024.05.01 10:00-05.01 14:30"
$ws_sheet1.Range("F7").Value = 102
$ws_sheet1.Range("G7").Value = 528
$ws_sheet1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=83863"
$ws_sheet1.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202404/oLIpAQh21712485244287.jpeg"
$ws_sheet1.Range("C8").Value = "北京·卡淘嘉年华·第三届球星卡交流会"
$ws_sheet1.Range("E8").Value = "2024.05.01 09:30-05.03 17:00"
$ws_sheet1.Range("F8").Value = 168
$ws_sheet1.Range("G8").Value = 85
$ws_sheet1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82072"
$ws_sheet1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202402/XOTabMFt1708929919204.jpeg"
$ws_sheet1.Range("C9").Value = "北京·原神x穹铁北京同人嘉年华7th"
$ws_sheet1.Range("D9").Value = "天辰东路7号 北京国家会议中心"
$ws_sheet1.Range("E9").Value = "2024.05.01 09:00-05.04 17:00"
$ws_sheet1.Range("F9").Value = 312
$ws_sheet1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84114"
$ws_sheet1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202404/55ApL1HY1712813894389.jpeg"
$ws_sheet1.Range("C10").Value = "北京·国乙同好嘉年华7th"
$ws_sheet1.Range("D10").Value = "北京国家会议中心 北京国家会议中心"
$ws_sheet1.Range("F10").Value = 420
$ws_sheet1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=82391"
$ws_sheet1.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202403/BGYIf9qe1709696198696.jpeg"
$ws_sheet1.Range("C11").Value = "北京·广播剧《宝石商人和钻石小姐》专场活动"
$ws_sheet1.Range("D11").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws_sheet1.Range("E11").Value = "2024.05.01 12:00-05.01 15:30"
$ws_sheet1.Range("F11").Value = 154
$ws_sheet1.Range("G11").Value = 288
$ws_sheet1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=82905"
$ws_sheet1.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202403/4RL1kiJi1710412443193.jpeg"
$ws_sheet1.Range("F12").Value = 199
$ws_sheet1.Range("F13").Value = 20
$ws_sheet1.Range("F14").Value = 460
$ws_sheet1.Range("F15").Value = 12254
$ws_sheet1.Range("F16").Value = 36
$ws_sheet1.Range("F18").Value = 308
$ws_sheet1.Range("F22").Value = 45
$ws_sheet1.Range("F23").Value = 156
$ws_sheet1.Range("F26").Value = 177
$ws_sheet1.Range("F27").Value = 161
$ws_sheet1.Range("F28").Value = 2734
$ws_sheet1.Range("F31").Value = 2103
$ws_sheet1.Range("F32").Value = 73
$ws_sheet1.Range("F34").Value = 2152
$ws_sheet1.Range("F35").Value = 1027
$ws_sheet1.Range("F36").Value = 4215
$ws_sheet1.Range("F37").Value = 3677
$ws_sheet1.Range("F38").Value = 613
$ws_sheet1.Range("F39").Value = 2629
$ws_sheet1.Range("F42").Value = 29
$ws_sheet1.Range("F43").Value = 116
$ws_sheet1.Range("F44").Value = 555
$ws_sheet1.Range("F45").Value = 71
$ws_sheet1.Range("F46").Value = 140
$ws_sheet1.Range("F47").Value = 238
$ws_sheet1.Range("F49").Value = 132

$ws_sheet2 = $wb.Worksheets.Item("演出")
$ws_sheet2.Range("F15").Value = 24
$ws_sheet2.Range("F17").Value = 32
$ws_sheet2.Range("F24").Value = 77

$ws_sheet4 = $wb.Worksheets.Item("全部类型")
$ws_sheet4.Range("F5").Value = 257
$ws_sheet4.Range("F6").Value = 9626
$ws_sheet4.Range("F7").Value = 642
$ws_sheet4.Range("F11").Value = 312
$ws_sheet4.Range("F12").Value = 420
$ws_sheet4.Range("F13").Value = 199
$ws_sheet4.Range("F14").Value = 20
$ws_sheet4.Range("F15").Value = 460
$ws_sheet4.Range("F16").Value = 12254
$ws_sheet4.Range("F17").Value = 36
$ws_sheet4.Range("F18").Value = 308
$ws_sheet4.Range("F22").Value = 156
$ws_sheet4.Range("F26").Value = 177
$ws_sheet4.Range("F27").Value = 161
$ws_sheet4.Range("F28").Value = 2734
$ws_sheet4.Range("F29").Value = 2103
$ws_sheet4.Range("F30").Value = 73
$ws_sheet4.Range("F31").Value = 2152
$ws_sheet4.Range("F32").Value = 1027
$ws_sheet4.Range("F36").Value = 4215
$ws_sheet4.Range("F37").Value = 3677
$ws_sheet4.Range("F38").Value = 613
$ws_sheet4.Range("F39").Value = 2629
$ws_sheet4.Range("F42").Value = 29
$ws_sheet4.Range("F43").Value = 116
$ws_sheet4.Range("F44").Value = 555
$ws_sheet4.Range("F45").Value = 71
$ws_sheet4.Range("F46").Value = 140
$ws_sheet4.Range("F47").Value = 238
$ws_sheet4.Range("F49").Value = 132

